$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Multi-line parameter strings (exact text, incl. embedded newlines) ----
$s34 = "learning_rate': 0.2, 'max_delta_step': 0, 'max_depth': 4, 'min_child_weight': 1, 'n_estimators': 800"
$s35 = "bootstrap': False,`n 'criterion': 'gini',`n 'max_depth': None,`n 'max_features': 10,`n 'min_samples_leaf': 10,`n 'min_samples_split': 2,`n 'n_estimators': 300"
$s36 = "{'class_weight': 'balanced',`n 'learning_rate': 0.01,`n 'max_depth': 10,`n 'min_child_samples': 200,`n 'n_estimators': 800,`n 'num_leaves': 24}"
$s37 = "{'class_weight': 'balanced',`n 'learning_rate': 0.01,`n 'max_depth': 8,`n 'min_child_samples': 400,`n 'n_estimators': 900,`n 'num_leaves': 24}"
$s38 = "bootstrap': False,`n 'criterion': 'gini',`n 'max_depth': 10,`n 'max_features': 1,`n 'min_samples_leaf': 1,`n 'min_samples_split': 3,`n 'n_estimators': 100"
$s39 = "class_weight='balanced', n_estimators=800, learning_rate=0.01, max_depth=8, min_child_samples=400, num_leaves=24"
$s40 = "class_weight='balanced', n_estimators=900, learning_rate=0.01, max_depth=8, min_child_samples=400, num_leaves=24"
$s41 = "class_weight='balanced', n_estimators=1000, learning_rate=0.01, max_depth=8, min_child_samples=400, num_leaves=24"

$dataset = "preprocessed_train_val_Mar13_0130pm_label_enc"
$lgbm = "LGBM"
$rf = "RandomForest"

# ---- Row 12 / 13: add Dataset column value (column E) ----
$ws.Range("E12").Value = $dataset
$ws.Range("E13").Value = $dataset

# ---- Row 14: finish the XGBClassifier row (B, E, I) ----
$ws.Range("B14").Value = $s34
$ws.Range("B13").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("E14").Value = $dataset
$ws.Range("I14").Value = 0.38369999999999999
$ws.Range("I13").Copy()
$ws.Range("I14").PasteSpecial(-4122)

# ---- Row 15 (new): RandomForest ----
$ws.Range("A15").Value = $rf
$ws.Range("B15").Value = $s35
$ws.Range("B9").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("E15").Value = $dataset
$ws.Range("I15").Value = 0.39945000000000003
$ws.Range("I13").Copy()
$ws.Range("I15").PasteSpecial(-4122)
$ws.Rows.Item(15).RowHeight = 101.5

# ---- Row 16 (new): LGBM ----
$ws.Range("A16").Value = $lgbm
$ws.Range("B16").Value = $s36
$ws.Range("B16").WrapText = $true
$ws.Range("E16").Value = $dataset
$ws.Range("H16").Value = 0.34752
$ws.Range("H16").NumberFormat = "0.00%"
$ws.Range("H16").Interior.Color = 49407
$ws.Rows.Item(16).RowHeight = 87

# ---- Row 17 (new): LGBM ----
$ws.Range("A17").Value = $lgbm
$ws.Range("B17").Value = $s37
$ws.Range("B17").WrapText = $true
$ws.Range("H17").Value = 0.34945999999999999
$ws.Range("H17").NumberFormat = "0.00%"
$ws.Range("H17").Interior.Color = 49407
$ws.Rows.Item(17).RowHeight = 87

# ---- Row 18 (new): RandomForest ----
$ws.Range("A18").Value = $rf
$ws.Range("B18").Value = $s38
$ws.Range("B9").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("H18").Value = 0.35060000000000002
$ws.Range("H18").NumberFormat = "0.00%"
$ws.Range("H18").Interior.Color = 49407
$ws.Rows.Item(18).RowHeight = 101.5

# ---- Row 19 (new): blank formatted cell ----
$ws.Range("B9").Copy()
$ws.Range("B19").PasteSpecial(-4122)

# ---- Row 20 (new): LGBM ----
$ws.Range("A20").Value = $lgbm
$ws.Range("B20").Value = $s39
$ws.Range("H20").Value = 0.34689999999999999
$ws.Range("F2").Copy()
$ws.Range("H20").PasteSpecial(-4122)

# ---- Row 21 (new): LGBM ----
$ws.Range("A21").Value = $lgbm
$ws.Range("B21").Value = $s40
$ws.Range("H21").Value = 0.34945999999999999
$ws.Range("F2").Copy()
$ws.Range("H21").PasteSpecial(-4122)

# ---- Row 22 (new) ----
$ws.Range("B22").Value = $s41
$ws.Range("H22").Value = 0.34660000000000002
$ws.Range("F2").Copy()
$ws.Range("H22").PasteSpecial(-4122)

# ---- Selection / view state ----
$ws.Range("E27").Select()
